$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($row, $col, $val)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue 2 4 "67.834.69"
Set-TextValue 2 5 "  +1.30%  "
Set-TextValue 3 4 "2.542.05"
Set-TextValue 3 5 "  +0.66%  "
Set-TextValue 4 5 "  -0.03%  "
Set-TextValue 5 4 "591.75"
Set-TextValue 5 5 "  +0.46%  "
Set-TextValue 6 4 "173.57"
Set-TextValue 6 5 "  +0.60%  "
Set-TextValue 7 5 "  -0.03%  "
Set-TextValue 8 5 "  -0.06%  "
Set-TextValue 9 4 "2.540.60"
Set-TextValue 9 5 "  +0.63%  "
Set-TextValue 10 4 "0.139"
Set-TextValue 10 5 "  +0.65%  "
Set-TextValue 11 5 "  +1.30%  "
Set-TextValue 12 5 "  -0.29%  "
Set-TextValue 14 2 "Binance-PegBSC-USD"
Set-TextValue 14 3 "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue 14 4 "2.54"
Set-TextValue 14 5 "  +154.62%  "
Set-TextValue 15 2 "Avalanche"
Set-TextValue 15 3 "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue 15 4 "26.56"
Set-TextValue 15 5 "  +0.24%  "
Set-TextValue 16 2 "ShibaInu"
Set-TextValue 16 3 "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue 16 4 "0.0000178"
Set-TextValue 16 5 "  +1.22%  "
Set-TextValue 17 2 "WrappedliquidstakedEther2.0"
Set-TextValue 17 3 "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue 17 4 "2.932.79"
Set-TextValue 17 5 "  -1.80%  "
Set-TextValue 18 2 "WrappedBTC"
Set-TextValue 18 3 "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue 18 4 "67.557.97"
Set-TextValue 18 5 "  +1.00%  "
Set-TextValue 19 2 "WrappedEther"
Set-TextValue 19 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue 19 4 "2.532.21"
Set-TextValue 19 5 "  +0.48%  "
Set-TextValue 20 2 "Chainlink"
Set-TextValue 20 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue 20 4 "11.82"
Set-TextValue 20 5 "  +4.46%  "
Set-TextValue 21 2 "Uniswap"
Set-TextValue 21 3 "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue 21 4 "7.97"
Set-TextValue 21 5 "  -1.56%  "
Set-TextValue 22 2 "BitcoinCash"
Set-TextValue 22 3 "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue 22 4 "370.97"
Set-TextValue 22 5 "  +4.54%  "
Set-TextValue 23 2 "Polkadot"
Set-TextValue 23 3 "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue 23 4 "4.16"
Set-TextValue 23 5 "  -0.31%  "
Set-TextValue 24 2 "NEARProtocol"
Set-TextValue 24 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue 24 4 "4.59"
Set-TextValue 24 5 "  -0.36%  "
Set-TextValue 25 2 "Litecoin"
Set-TextValue 25 3 "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue 25 4 "71.61"
Set-TextValue 25 5 "  +2.87%  "
Set-TextValue 26 2 "Dai"
Set-TextValue 26 3 "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue 26 4 "1.00"
Set-TextValue 26 5 "  -0.04%  "
Set-TextValue 27 2 "SuiNetwork"
Set-TextValue 27 3 "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue 27 4 "1.93"
Set-TextValue 27 5 "  -2.73%  "
Set-TextValue 28 2 "Aptos"
Set-TextValue 28 3 "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue 28 4 "10.00"
Set-TextValue 28 5 "  +0.58%  "
Set-TextValue 29 2 "WrappedeETH"
Set-TextValue 29 3 "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue 29 4 "2.665.24"
Set-TextValue 29 5 "  +0.43%  "
Set-TextValue 30 2 "PEPE"
Set-TextValue 30 3 "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue 30 4 "0.0₃0970"
Set-TextValue 30 5 "  -0.33%  "
Set-TextValue 31 2 "InternetComputer(DFINITY)"
Set-TextValue 31 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue 31 4 "8.51"
Set-TextValue 31 5 "  +4.84%  "
Set-TextValue 32 2 "Bittensor"
Set-TextValue 32 3 "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue 32 4 "539.24"
Set-TextValue 32 5 "  +1.53%  "
Set-TextValue 33 2 "Fetch.AI"
Set-TextValue 33 3 "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue 33 4 "1.32"
Set-TextValue 33 5 "  -0.12%  "
Set-TextValue 34 2 "PancakeSwap"
Set-TextValue 34 3 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue 34 4 "1.88"
Set-TextValue 34 5 "  +1.73%  "
Set-TextValue 35 2 "Kaspa"
Set-TextValue 35 3 "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue 35 4 "0.129"
Set-TextValue 35 5 "  -1.50%  "
Set-TextValue 36 2 "FirstDigitalUSD"
Set-TextValue 36 3 "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue 36 4 "0.999"
Set-TextValue 36 5 "  -0.06%  "
Set-TextValue 37 2 "Monero"
Set-TextValue 37 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 37 4 "159.13"
Set-TextValue 37 5 "  +1.40%  "
Set-TextValue 38 2 "ImmutableX"
Set-TextValue 38 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue 38 4 "1.44"
Set-TextValue 38 5 "  -1.18%  "
Set-TextValue 39 2 "EthereumClassic"
Set-TextValue 39 3 "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue 39 4 "19.22"
Set-TextValue 39 5 "  +3.44%  "
Set-TextValue 40 2 "WhiteBITCoin"
Set-TextValue 40 3 "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue 40 4 "18.63"
Set-TextValue 40 5 "  +1.07%  "
Set-TextValue 41 5 "  +0.97%  "
Set-TextValue 42 2 "PolygonEcosystemToken"
Set-TextValue 42 3 "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue 42 4 "0.352"
Set-TextValue 42 5 "  -0.33%  "
Set-TextValue 43 2 "Stacks"
Set-TextValue 43 3 "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue 43 4 "1.78"
Set-TextValue 43 5 "  -0.38%  "
Set-TextValue 44 2 "dogwifhat"
Set-TextValue 44 3 "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue 44 4 "2.58"
Set-TextValue 44 5 "  +3.94%  "
Set-TextValue 45 2 "USDe"
Set-TextValue 45 3 "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue 45 4 "1.00"
Set-TextValue 45 5 "  +0.02%  "
Set-TextValue 46 2 "OKB"
Set-TextValue 46 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue 46 4 "39.18"
Set-TextValue 46 5 "  -1.27%  "
Set-TextValue 47 2 "BabyDogeCoin"
Set-TextValue 47 3 "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue 47 4 "0.0₆0292"
Set-TextValue 47 5 "  +5.32%  "
Set-TextValue 48 2 "Aave"
Set-TextValue 48 3 "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue 48 4 "147.63"
Set-TextValue 48 5 "  -0.95%  "
Set-TextValue 49 2 "ARBITRUM"
Set-TextValue 49 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue 49 4 "0.555"
Set-TextValue 49 5 "  +0.05%  "
Set-TextValue 50 2 "Filecoin"
Set-TextValue 50 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue 50 4 "3.72"
Set-TextValue 50 5 "  +0.94%  "
Set-TextValue 51 2 "Optimism"
Set-TextValue 51 3 "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
Set-TextValue 51 4 "1.72"
Set-TextValue 51 5 "  +2.26%  "
